$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new row at the top; this pushes every existing
#    municipio row down by one (row 1 -> row 2, row 2 -> row 3, ...).
$ws.Range("A1").EntireRow.Insert()

# 2) Populate the new header row.
$ws.Range("A1").Value = "MUNICIPIO"
$ws.Range("B1").Value = "CASOS"
$ws.Range("C1").Value = [char]0x00D3 + "BITOS"

# 3) Style the header row: thin box border, bold font, centered
#    horizontally, top-aligned vertically.
$headerRng = $ws.Range("A1:C1")
$headerRng.Borders.LineStyle = 1
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# 4) Append the two brand-new rows at the bottom of the table
#    (after the insert, the old last row "votuporanga" now sits on
#    row 155, so the new data goes on rows 156-157).
$newRows = New-Object 'object[,]' 2,3
$newRows[0,0] = "outros estados"
$newRows[0,1] = 49
$newRows[0,2] = ""
$newRows[1,0] = "outros paises"
$newRows[1,1] = 40
$newRows[1,2] = ""
$ws.Range("A156:C157").Value = $newRows
